$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8549122214317322
$ws.Range("B1").Value = 1.361811995506287
$ws.Range("C1").Value = 2.745779991149902
$ws.Range("D1").Value = 3.878943920135498
$ws.Range("E1").Value = 1.876351356506348
